$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename columns: year_quarter -> cal_year_quarter, year_month -> cal_year_month
$ws.Range("A18").Value = "cal_year_quarter"
$ws.Range("A19").Value = "cal_year_month"

# Update the active selection on the sheet to A19
[void]$ws.Range("A19").Select()
